# Applies the "Add files via upload" revision to the HDSHM congress schedule
# workbook: updates the poster-session summary texts in row 24 (authors list,
# Croatian poster list, English poster list) and appends a new "Lunch" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 24 (poster presentation row): refresh the author list plus the
#    Croatian/English poster-list summaries with the updated wording (adds
#    Anastazija Ramljak & Anita Krnjak, fixes "Stefanec" -> "Stefanes", etc.)
# ---------------------------------------------------------------------------

$authors = @"
Zvonimir Banoža, Matija Radobuljac, Sanja Pavić Jelečki, Suzana Palatinuš, Anastazija Ramljak, Anita Krnjak, Ivana Kralj, Saša Balija, Silvia Tisaj Pigac, Ivana Sklepić Klobučarić, Tina Kresonja, Milena Škvorc, Suzana Palatinuš, Barbara Samvik, Kristijan Štefanes, Nikola Čopor, Ivica Bračko, Ramon Tumbas, Emil Kralj, Dragica Svetličić
"@

$posterHr = @"
P—1    Anafilaktički šok,
Zvonimir Banoža, Matija Radobuljac
P—2    Disekcija aorte
Sanja Pavić Jelečki
P—3    Edukacija zaposlenika Zavoda za hitnu medicinu
Krapinsko-zagorske županije
Suzana Palatinuš
P—4    Inhalacijska analgezija u hitnoj medicini: primjena Penthroxa
Anastazija Ramljak, Anita Krnjak
P—5    Inovativni pristupi u edukaciji hitnih stanja
Ivana Kralj, Saša Balija
P—6    Multidisciplinarni pristup u zbrinjavanju trudnice s abrupcijom placente od HMP-a, rađaone do operacijske sale
Silvia Tisaj Pigac, Ivana Sklepić Klobučarić, Tina Kresonja, Milena Škvorc
P—7    Reanimacija u izvanbolničkim uvjetima
Suzana Palatinuš
P—8    Stavovi građana o oživljavanju osoba van bolnice i DNR (Do Not Resuscitate) - “Ne oživljavaj” obrascu
Barbara Samvik, Kristijan Štefanes, Nikola Čopor, Ivica Bračko
P—9    Važnost prepoznavanja simptoma i znakova kardiogenog šoka kod pacijenata u OHBP-u s osvrtom na slučaj pacijenta s postinfarktnom rupturom septuma
Ramon Tumbas, Emil Kralj, Dragica Svetličić
"@

$posterEn = @"
P—1  Anaphylactic shock,
Zvonimir Banoža, Matija Radobuljac
P—2  Aortic dissection
Sanja Pavić Jelečki
P—3  Education of employees of the Institute of Emergency Medicine of Krapina - Zagorje County
Suzana Palatinuš
P—4  Inhalation Analgesia in Emergency Medicine: Application of Penthrox
Anastazija Ramljak, Anita Krnjak
P—5  Innovative Approaches in Emergency Education
Ivana Kralj, Saša Balija
P—6   Multidisciplinary Approach in Managing a Pregnant Woman with Placental Abruption – From EMS to Delivery Room and Operating Theatre
Silvia Tisaj Pigac, Ivana Sklepić Klobučarić, Tina Kresonja, Milena Škvorc
P—7 Resuscitation in out-of-hospital conditions
Suzana Palatinuš
P—8 Stavovi građana o oživljavanju osoba van bolnice i DNR (Do Not Resuscitate) - “Ne oživljavaj” obrascu
Barbara Samvik, Kristijan Štefanes, Nikola Čopor, Ivica Bračko
P—9 The importance of recognizing the symptoms and signs of cardiogenic shock in patients in the ED with reference to the case of a patient with postinfarction septal rupture
Ramon Tumbas, Emil Kralj, Dragica Svetličić
"@

$ws.Cells.Item(24, 5).Value = $authors
$ws.Cells.Item(24, 11).Value = $posterHr
$ws.Cells.Item(24, 12).Value = $posterEn

# ---------------------------------------------------------------------------
# 2. New row 27: lunch break entry (mirrors the layout of the other
#    break-type rows such as row 25/26 - col A holds the generic "break"
#    type marker, B the time span, C/D the HR/EN titles, J the hall column).
# ---------------------------------------------------------------------------

$ws.Cells.Item(27, 1).Value = $ws.Cells.Item(26, 1).Value()
$ws.Cells.Item(27, 2).Value = "14:00 - 16:00"
$ws.Cells.Item(27, 3).Value = "Ručak"
$ws.Cells.Item(27, 4).Value = "Lunch"
$ws.Cells.Item(27, 10).Value = $ws.Cells.Item(13, 10).Value()
